$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns C (fcs), E (MYDIR), F (need_to_buy) for rows 2-15.
# Column B (TISG) and D (buy) are left untouched.

$updates = @{
    2  = @{ C = 14003.0577545912; E = 9891.95928647022; F = 12.4648767108917 }
    3  = @{ C = 14055.2447855365; E = 9909.46166717098; F = 495.368602196143 }
    4  = @{ C = 12382.3344583256; E = 9665.60310024208; F = 415.503231606987 }
    5  = @{ C = 8753.91106920553; E = 8937.0973935836;  F = 233.96451928288  }
    6  = @{ C = 8925.99306536615; E = 8617.69090347708; F = 227.825998701801 }
    7  = @{ C = 13457.5083027979; E = 9107.15889740035; F = 437.033633341592 }
    8  = @{ C = 13045.216970397;                        F = 419.853348513342 }
    9  = @{ C = 12492.1984726227;                       F = 396.810911106078 }
    10 = @{ C = 12086.9072756475;                       F = 379.92377789878  }
    11 = @{ C = 11225.0234085016;                       F = 344.010280077118 }
    12 = @{ C = 7558.89536473484;                       F = 173.971156765581 }
    13 = @{ C = 7500.33258681163;                       F = 171.013712026297 }
    14 = @{ C = 11576.4454605222;                       F = 349.299287608749 }
    15 = @{ C = 11833.9242362269;                       F = 360.027569929778 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("C$row").Value = $vals.C
    if ($vals.ContainsKey("E")) {
        $ws.Range("E$row").Value = $vals.E
    }
    $ws.Range("F$row").Value = $vals.F
}
